$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted ahead of the existing series
# (which pushes every following row down by one and appends a duplicate
# of the former last row at the new bottom, row 307).
$ws.Rows.Item(247).Insert()

$ws.Range("A247").Value = 8
$ws.Range("B247").Value = "Terminal La Palmera de La Serena"
$ws.Range("C247").Value = "Coquimbo"
$ws.Range("D247").Value = 44798
$ws.Range("E247").Value = 4
$ws.Range("F247").Value = 100112012
$ws.Range("G247").Value = "Espinaca"
$ws.Range("H247").Value = "Sin especificar"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 2000
$ws.Range("K247").Value = 500
$ws.Range("L247").Value = 550
$ws.Range("M247").Value = 525
$ws.Range("N247").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O247").Value = "Provincia del Elquí"
$ws.Range("P247").Value = 1050
$ws.Range("Q247").Value = 0.5
$ws.Range("R247").Value = "Hortaliza"
